$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.092.28"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.66%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.175.70"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.64%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.11%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.83"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.13%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "150.99"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.01%  "

# Row 7
$ws.Range("E7").Value = "  -0.01%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.173.30"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.65%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.534"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.79%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.163"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.08%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.22"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.86%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.505"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.92%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000270"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +16.87%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "38.01"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.91%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.689.57"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.54%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.111.03"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.65%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.168.64"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.52%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.18"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +5.21%  "

# Row 19
$ws.Range("E19").Value = "  +1.38%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "513.27"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +7.29%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.87"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.89%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.733"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.66%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "15.33"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.62%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.82"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.05%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.43"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.36%  "

# Row 26
$ws.Range("E26").Value = "  +0.12%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.07"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +11.07%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.94"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.19%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.19"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +7.12%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "28.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.19%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.80"
$ws.Range("D31").Style = "Normal"

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.999"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.03%  "

# Row 33
$ws.Range("E33").Value = "  +5.22%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.33"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +9.09%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.61"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.56%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "55.79"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.99%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0899"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +10.05%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "477.66"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.71%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.16"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +11.26%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0422"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.16%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.66"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.17%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.064.32"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.26%  "

# Row 43
$ws.Range("E43").Value = "  +3.03%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.287"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.98%  "

# Row 45
$ws.Range("E45").Value = "  +6.35%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "29.18"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.37%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0₃0618"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +19.00%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.00"
$ws.Range("D48").Style = "Normal"

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.116"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.00%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.27"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +7.16%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "120.63"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.37%  "
